# The workbook's single worksheet is renamed from "Time" to "C.Size" and the
# active selection is moved from G3:G17 (anchored at G3) to the single cell J11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "C.Size"

$ws.Range("J11").Select()
